$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "username"
$ws.Range("A3").Value = "password"

$ws.Range("B2").Value = "anas123"
$ws.Range("B3").Value = "Zehra@2016"

$ws.Range("A4").Value = "userAccountName"
$ws.Range("B4").Value = "anas"

[void]$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Zehra@2016")

[void]$ws.Range("B7").Select()
